$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "Bad" (red fill / red text) cell style to highlight the
# "blocks stats" rows that needed review: rows 12, 13, 14 and 18 (A:H).
$ws.Range("A12:H12").Style = "Bad"
$ws.Range("A13:H13").Style = "Bad"
$ws.Range("A14:H14").Style = "Bad"
$ws.Range("A18:H18").Style = "Bad"

# Remove the leftover duplicated "Wither Ske" / "KB 3" row data in A24:B24.
$ws.Range("A24:B24").ClearContents()

# Leave the selection where the author finished editing.
$ws.Range("D29").Select()
